# Fruta / hortaliza, semanal
# A new weekly price record is inserted at row 221 (before the existing
# "Frutilla" records), pushing all subsequent rows (old 221-303) down by
# one row to (222-304).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 221; this shifts rows 221:303 down to
# 222:304 and picks up the formatting (including the date-cell style) of
# the row above, matching the existing layout.
$ws.Rows.Item(221).Insert()

# Populate the newly inserted row 221 with the new data record.
$ws.Range("A221").Value2 = 10
$ws.Range("B221").Value2 = "Vega Modelo de Temuco"
$ws.Range("C221").Value2 = "La Araucanía"
$ws.Range("D221").Value2 = 44917
$ws.Range("E221").Value2 = 9
$ws.Range("F221").Value2 = "Fruta"
$ws.Range("G221").Value2 = 100101
$ws.Range("H221").Value2 = "Berries"
$ws.Range("I221").Value2 = 100112025
$ws.Range("J221").Value2 = "Frutilla"
$ws.Range("K221").Value2 = "Sin especificar"
$ws.Range("L221").Value2 = "Primera"
$ws.Range("M221").Value2 = 300
$ws.Range("N221").Value2 = 8500
$ws.Range("O221").Value2 = 8500
$ws.Range("P221").Value2 = 8500
$ws.Range("Q221").Value2 = "$/caja 7 kilos"
$ws.Range("R221").Value2 = "Región de La Araucanía"
$ws.Range("S221").Value2 = 1214
$ws.Range("T221").Value2 = 7
